$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 8: userMail = jamal@mytestmail.net (existing value reused), userPass = salah.Passwd (new value)
$ws.Range("A8").Value = "jamal@mytestmail.net"
$ws.Range("B8").Value = "salah.Passwd"

# Move active selection to B9 (as in the post-edit file)
$ws.Range("B9").Select()
